# Insert a new data row before the current row 167, shifting the existing
# rows 167-193 down to 168-194, and populate the newly inserted row with
# the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 167:193 down by inserting a new blank row at 167.
$ws.Rows("167:167").Insert()

# Populate the new row 167 with the new record.
$ws.Cells.Item(167, 1).Value  = 3
$ws.Cells.Item(167, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(167, 3).Value  = "Coquimbo"
$ws.Cells.Item(167, 4).Value  = 44711
$ws.Cells.Item(167, 5).Value  = 5
$ws.Cells.Item(167, 6).Value  = 100112030
$ws.Cells.Item(167, 7).Value  = "Poroto granado"
$ws.Cells.Item(167, 8).Value  = "Sin especificar"
$ws.Cells.Item(167, 9).Value  = "Primera"
$ws.Cells.Item(167, 10).Value = 38
$ws.Cells.Item(167, 11).Value = 21000
$ws.Cells.Item(167, 12).Value = 21000
$ws.Cells.Item(167, 13).Value = 21000
$ws.Cells.Item(167, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(167, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(167, 16).Value = 840
$ws.Cells.Item(167, 17).Value = 25
$ws.Cells.Item(167, 18).Value = "Hortaliza"
